$d = $word.ActiveDocument

# ------------------------------------------------------------------
# This edit changes "nodes" -> "blocks" in two places inside the
# "For thread 1, ..." paragraph. In both places, Word ends up
# splitting the paragraph's run so that the word "blocks" sits in its
# own <w:r>. The first occurrence additionally carries the document's
# "_GoBack" bookmark (the automatic bookmark Word drops at the most
# recent edit position) immediately after the new word; since a
# bookmark name must be unique, adding a new "_GoBack" bookmark here
# automatically removes the old one that used to sit at the end of
# the "Thread 3" paragraph.
#
# To reproduce the run split cleanly (without leaving a stray
# "blocks" run with an unnecessary xml:space="preserve"), we drop a
# temporary bookmark at the boundary *before* replacing the word, and
# remove the temporary bookmark again once the split has been
# created.
# ------------------------------------------------------------------

# --- First occurrence: "...number of nodes in the free list ..." ---
$r = $d.Content
$r.Find.Execute("nodes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$wordStart = $r.Start
$wordEnd = $r.End

# Mark the boundary before the word so the replacement below splits
# into its own run instead of re-merging with the preceding text.
$splitBefore = $d.Range($wordStart, $wordStart)
$d.Bookmarks.Add("_tmpSplitBefore1", $splitBefore) | Out-Null

# Mark the boundary after the word as "_GoBack" (this both forces the
# trailing split and places the real bookmark exactly where the diff
# expects it; it also evicts the old "_GoBack" bookmark elsewhere in
# the document because bookmark names are unique).
$splitAfter = $d.Range($wordEnd, $wordEnd)
$d.Bookmarks.Add("_GoBack", $splitAfter) | Out-Null

$target = $d.Range($wordStart, $wordEnd)
$target.Text = "blocks"

$d.Bookmarks("_tmpSplitBefore1").Delete()

# --- Second occurrence: "...if there are no nodes in the free list..." ---
$r2 = $d.Range($target.End, $d.Content.End)
$r2.Find.Execute("nodes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$wordStart2 = $r2.Start
$wordEnd2 = $r2.End

$splitBefore2 = $d.Range($wordStart2, $wordStart2)
$d.Bookmarks.Add("_tmpSplitBefore2", $splitBefore2) | Out-Null

$splitAfter2 = $d.Range($wordEnd2, $wordEnd2)
$d.Bookmarks.Add("_tmpSplitAfter2", $splitAfter2) | Out-Null

$target2 = $d.Range($wordStart2, $wordEnd2)
$target2.Text = "blocks"

$d.Bookmarks("_tmpSplitBefore2").Delete()
$d.Bookmarks("_tmpSplitAfter2").Delete()
